{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The document ends with two empty paragraphs before the sectPr.\n// Insert the new \"CLAUDE\" label paragraph and the Claude API key\n// paragraph right before the very last (empty) paragraph, i.e.\n// immediately after the second-to-last (empty) paragraph.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nconst labelParagraph = lastParagraph.insertParagraph(\"CLAUDE\", \"Before\");\nlabelParagraph.insertParagraph(\n  \"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\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The document ends with two empty paragraphs right before the section\n# break. We insert the new \"CLAUDE\" label paragraph and the Claude API\n# key paragraph between them, i.e. right after the second-to-last\n# (empty) paragraph and before the very last (empty) paragraph.\n$count = $d.Paragraphs.Count\n$secondToLast = $d.Paragraphs.Item($count - 1)\n\n$r = $secondToLast.Range\n$r.Collapse(0)          # wdCollapseEnd\n$r.InsertParagraphAfter()\n\n$labelPara = $d.Paragraphs.Item($count)\n$labelPara.Range.InsertAfter(\"CLAUDE\")\n\n$r2 = $labelPara.Range\n$r2.Collapse(0)          # wdCollapseEnd\n$r2.InsertParagraphAfter()\n\n$keyPara = $d.Paragraphs.Item($count + 1)\n$keyPara.Range.InsertAfter(\"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\")\n"}
